$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 338.73077
$ws.Range("I19").Value = 274.41666
$ws.Range("J19").Value = 393.85715
$ws.Range("K19").Value = 274.41666
$ws.Range("L19").Value = 393.85715
$ws.Range("M19").Value = -99.41665999999998
$ws.Range("N19").Value = -743.85715
# Row 51
$ws.Range("H51").Value = 3032.12
$ws.Range("I51").Value = 2934
$ws.Range("J51").Value = 3045.5
$ws.Range("K51").Value = 2934
$ws.Range("L51").Value = 3045.5
$ws.Range("M51").Value = -2450
$ws.Range("N51").Value = -4013.5
# Row 129
$ws.Range("H129").Value = 1979.6666
$ws.Range("I129").Value = 300
$ws.Range("J129").Value = 2189.625
$ws.Range("K129").Value = 900
$ws.Range("L129").Value = 6568.875
$ws.Range("M129").Value = 4100
$ws.Range("N129").Value = -16568.875

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 5402.087
$ws.Range("I74").Value = 1032.6316
$ws.Range("J74").Value = 26157
$ws.Range("K74").Value = 1032.6316
$ws.Range("L74").Value = 26157
$ws.Range("M74").Value = -158.6315999999999
$ws.Range("N74").Value = -27905
# Row 77
$ws.Range("H77").Value = 5402.087
$ws.Range("I77").Value = 1032.6316
$ws.Range("J77").Value = 26157
$ws.Range("K77").Value = 5163.157999999999
$ws.Range("L77").Value = 130785
$ws.Range("M77").Value = -795.1579999999994
$ws.Range("N77").Value = -139521
# Row 134
$ws.Range("H134").Value = 39062.332
$ws.Range("J134").Value = 39062.332
$ws.Range("L134").Value = 39062.332
$ws.Range("N134").Value = -49202.332
# Row 141
$ws.Range("H141").Value = 29429
$ws.Range("J141").Value = 29429
$ws.Range("L141").Value = 29429
$ws.Range("N141").Value = -39789

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 15428.454
$ws.Range("I82").Value = 7745.8887
$ws.Range("J82").Value = 50000
$ws.Range("K82").Value = 7745.8887
$ws.Range("L82").Value = 50000
$ws.Range("M82").Value = -7362.8887
$ws.Range("N82").Value = -50766
# Row 85
$ws.Range("H85").Value = 15428.454
$ws.Range("I85").Value = 7745.8887
$ws.Range("J85").Value = 50000
$ws.Range("K85").Value = 7745.8887
$ws.Range("L85").Value = 50000
$ws.Range("M85").Value = -6419.8887
$ws.Range("N85").Value = -52652

$ws = $wb.Worksheets.Item("CRP")
# Row 17
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 1000
$ws.Range("M17").Value = -826
# Row 31
$ws.Range("H31").Value = 21715.715
$ws.Range("I31").Value = 1545.5
$ws.Range("J31").Value = 41885.93
$ws.Range("K31").Value = 1545.5
$ws.Range("L31").Value = 41885.93
$ws.Range("M31").Value = -1250.5
$ws.Range("N31").Value = -42475.93
# Row 34
$ws.Range("H34").Value = 21715.715
$ws.Range("I34").Value = 1545.5
$ws.Range("J34").Value = 41885.93
$ws.Range("K34").Value = 1545.5
$ws.Range("L34").Value = 41885.93
$ws.Range("M34").Value = -1343.5
$ws.Range("N34").Value = -42289.93
# Row 41
$ws.Range("H41").Value = 5500
# Row 50
$ws.Range("H50").Value = 24850
$ws.Range("J50").Value = 24850
$ws.Range("L50").Value = 24850
$ws.Range("N50").Value = -26100
# Row 51
$ws.Range("H51").Value = 14199.8
$ws.Range("I51").Value = 11499.75
$ws.Range("J51").Value = 25000
$ws.Range("K51").Value = 11499.75
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = -10763.75
$ws.Range("N51").Value = -26472
# Row 58
$ws.Range("H58").Value = 788.4545000000001
$ws.Range("I58").Value = 644.1429000000001
$ws.Range("J58").Value = 1041
$ws.Range("K58").Value = 644.1429000000001
$ws.Range("L58").Value = 1041
$ws.Range("M58").Value = -441.1429000000001
$ws.Range("N58").Value = -1447
# Row 60
$ws.Range("H60").Value = 10774.28
$ws.Range("I60").Value = 5833.3335
$ws.Range("J60").Value = 11448.046
$ws.Range("K60").Value = 5833.3335
$ws.Range("L60").Value = 11448.046
$ws.Range("M60").Value = -5322.3335
$ws.Range("N60").Value = -12470.046
# Row 61
$ws.Range("H61").Value = 14199.8
$ws.Range("I61").Value = 11499.75
$ws.Range("J61").Value = 25000
$ws.Range("K61").Value = 11499.75
$ws.Range("L61").Value = 25000
$ws.Range("M61").Value = -11151.75
$ws.Range("N61").Value = -25696
# Row 136
$ws.Range("H136").Value = 788.4545000000001
$ws.Range("I136").Value = 644.1429000000001
$ws.Range("J136").Value = 1041
$ws.Range("K136").Value = 1932.4287
$ws.Range("L136").Value = 3123
$ws.Range("M136").Value = 617.5712999999998
$ws.Range("N136").Value = -8223

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 28.516129
$ws.Range("J12").Value = 39.363636
$ws.Range("L12").Value = 118.090908
$ws.Range("N12").Value = -464.090908
# Row 80
$ws.Range("H80").Value = 2268.2856
$ws.Range("I80").Value = 2068.4
$ws.Range("J80").Value = 2450
$ws.Range("K80").Value = 6205.200000000001
$ws.Range("L80").Value = 7350
$ws.Range("M80").Value = -5269.200000000001
$ws.Range("N80").Value = -9222
# Row 81
$ws.Range("H81").Value = 1068.8572
$ws.Range("I81").Value = 312.25
$ws.Range("J81").Value = 1371.5
$ws.Range("K81").Value = 936.75
$ws.Range("L81").Value = 4114.5
$ws.Range("M81").Value = 186.25
$ws.Range("N81").Value = -6360.5
# Row 82
$ws.Range("H82").Value = 13333
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 13333
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 39999
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -40811
# Row 83
$ws.Range("H83").Value = 2268.2856
$ws.Range("I83").Value = 2068.4
$ws.Range("J83").Value = 2450
$ws.Range("K83").Value = 18615.6
$ws.Range("L83").Value = 22050
$ws.Range("M83").Value = -13935.6
$ws.Range("N83").Value = -31410
# Row 84
$ws.Range("H84").Value = 1068.8572
$ws.Range("I84").Value = 312.25
$ws.Range("J84").Value = 1371.5
$ws.Range("K84").Value = 2810.25
$ws.Range("L84").Value = 12343.5
$ws.Range("M84").Value = 2805.75
$ws.Range("N84").Value = -23575.5
# Row 85
$ws.Range("H85").Value = 13333
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 13333
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 39999
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -42807
# Row 131
$ws.Range("H131").Value = 704.3226
$ws.Range("I131").Value = 340.18918
$ws.Range("J131").Value = 944.9107
$ws.Range("K131").Value = 1020.56754
$ws.Range("L131").Value = 2834.7321
$ws.Range("M131").Value = 4019.43246
$ws.Range("N131").Value = -12914.7321

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1078.75
$ws.Range("I122").Value = 1078
$ws.Range("J122").Value = 1080
$ws.Range("K122").Value = 3234
$ws.Range("L122").Value = 3240
$ws.Range("M122").Value = -784
$ws.Range("N122").Value = -8140
# Row 127
$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
# Row 129
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
# Row 130
$ws.Range("H130").Value = 47890
$ws.Range("J130").Value = 47890
$ws.Range("L130").Value = 47890
$ws.Range("N130").Value = -57930
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 50006.285
$ws.Range("I122").Value = 79068.62
$ws.Range("J122").Value = 2780
$ws.Range("K122").Value = 237205.86
$ws.Range("L122").Value = 8340
$ws.Range("M122").Value = -234755.86
$ws.Range("N122").Value = -13240
# Row 132
$ws.Range("H132").Value = 4515.4165
$ws.Range("I132").Value = 5697.478
$ws.Range("J132").Value = 2424.077
$ws.Range("K132").Value = 17092.434
$ws.Range("L132").Value = 7272.231000000001
$ws.Range("M132").Value = -14562.434
$ws.Range("N132").Value = -12332.231

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 12501191
$ws.Range("I122").Value = 25000956
$ws.Range("J122").Value = 1426.25
$ws.Range("K122").Value = 75002868
$ws.Range("L122").Value = 4278.75
$ws.Range("M122").Value = -75000418
$ws.Range("N122").Value = -9178.75
# Row 126
$ws.Range("H126").Value = 766.871
$ws.Range("I126").Value = 643.8182
$ws.Range("K126").Value = 1931.4546
$ws.Range("M126").Value = 538.5454
